$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.377.82"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.942.38"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'242.04"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("E6").Value = "  -2.97%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'56.85"
$ws.Range("E8").Value = "  -3.96%  "
$ws.Range("E9").Value = "  -4.19%  "
$ws.Range("E10").Value = "  +4.39%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "2.228.29"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").Value = "'0.813"
$ws.Range("E13").Value = "  -5.98%  "
$ws.Range("D14").Value = "'13.45"
$ws.Range("E14").Value = "  -4.27%  "
$ws.Range("D15").Value = "'20.96"
$ws.Range("E15").Value = "  -12.31%  "
$ws.Range("D16").Value = "'5.15"
$ws.Range("E16").Value = "  -6.12%  "
$ws.Range("D17").Value = "1.951.97"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("D18").Value = "36.313.85"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "0.0₃0871"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "'69.26"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").Value = "'228.09"
$ws.Range("E21").Value = "  -2.81%  "
$ws.Range("D22").Value = "'5.00"
$ws.Range("E22").Value = "  -6.31%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -9.01%  "
$ws.Range("E25").Value = "  -1.64%  "
$ws.Range("D26").Value = "'9.20"
$ws.Range("E26").Value = "  -10.27%  "
$ws.Range("D27").Value = "'161.06"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("E28").Value = "  +3.37%  "
$ws.Range("D29").Value = "'19.22"
$ws.Range("E29").Value = "  -3.57%  "
$ws.Range("D30").Value = "'0.117"
$ws.Range("E30").Value = "  -2.64%  "
$ws.Range("E31").Value = "  -6.58%  "
$ws.Range("D32").Value = "'4.61"
$ws.Range("E32").Value = "  -6.45%  "
$ws.Range("D33").Value = "'0.0633"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("E34").Value = "  -4.62%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("E36").Value = "  -3.32%  "
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("E38").Value = "  -6.24%  "
$ws.Range("D39").Value = "'3.00"
$ws.Range("E39").Value = "  -3.06%  "
$ws.Range("D40").Value = "'0.0966"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("E43").Value = "  -7.29%  "
$ws.Range("D44").Value = "'15.61"
$ws.Range("E44").Value = "  -4.29%  "
$ws.Range("D45").Value = "1.341.69"
$ws.Range("E45").Value = "  -3.06%  "
$ws.Range("E46").Value = "  -6.99%  "
$ws.Range("D47").Value = "'86.85"
$ws.Range("E47").Value = "  -6.61%  "
$ws.Range("D48").Value = "'7.09"
$ws.Range("E48").Value = "  -6.46%  "
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("D50").Value = "'44.16"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").Value = "2.119.49"
$ws.Range("E51").Value = "  -2.42%  "
